$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.324.57"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.097.20"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.91"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.17"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.093.99"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").Value = "  -6.20%  "
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.20"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.596.84"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.341.87"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.097.25"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.56"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.38"
$ws.Range("E23").Value = "  +12.45%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.57"
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.38"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.80"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.03"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.70"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0410"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.99"
$ws.Range("E38").Value = "  +15.68%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "446.83"
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0816"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.967.97"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  -5.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.30"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.36"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("E50").Value = "  +0.87%  "
